# Apply cryptos list update (Fri Jan 12 18:27:22 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '43.659.67'
$ws.Cells.Item(2, 5).Value = '  -5.92%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.600.87'
$ws.Cells.Item(3, 5).Value = '  +0.17%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.00'
$ws.Cells.Item(4, 4).NumberFormat = "General"
$ws.Cells.Item(4, 5).Value = '  +0.00%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '301.24'
$ws.Cells.Item(5, 4).NumberFormat = "General"
$ws.Cells.Item(5, 5).Value = '  -2.26%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '96.46'
$ws.Cells.Item(6, 4).NumberFormat = "General"
$ws.Cells.Item(6, 5).Value = '  -3.60%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.580'
$ws.Cells.Item(7, 4).NumberFormat = "General"
$ws.Cells.Item(7, 5).Value = '  -4.03%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.08%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -3.22%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '37.11'
$ws.Cells.Item(10, 4).NumberFormat = "General"
$ws.Cells.Item(10, 5).Value = '  -5.54%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -3.42%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '7.86'
$ws.Cells.Item(12, 4).NumberFormat = "General"
$ws.Cells.Item(12, 5).Value = '  -4.02%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '2.993.70'
$ws.Cells.Item(13, 5).Value = '  +0.08%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  +1.13%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '2.593.35'
$ws.Cells.Item(15, 5).Value = '  +0.01%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  -3.11%  '

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '14.41'
$ws.Cells.Item(17, 4).NumberFormat = "General"
$ws.Cells.Item(17, 5).Value = '  -4.27%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '43.713.23'
$ws.Cells.Item(18, 5).Value = '  -6.05%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '6.69'
$ws.Cells.Item(19, 4).NumberFormat = "General"
$ws.Cells.Item(19, 5).Value = '  -1.17%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '0.0₃0982'
$ws.Cells.Item(20, 5).Value = '  -3.43%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '12.42'
$ws.Cells.Item(21, 4).NumberFormat = "General"
$ws.Cells.Item(21, 5).Value = '  -4.80%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '73.42'
$ws.Cells.Item(22, 4).NumberFormat = "General"
$ws.Cells.Item(22, 5).Value = '  +2.38%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '267.47'
$ws.Cells.Item(23, 4).NumberFormat = "General"
$ws.Cells.Item(23, 5).Value = '  -3.56%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +1.92%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.94'
$ws.Cells.Item(25, 4).NumberFormat = "General"
$ws.Cells.Item(25, 5).Value = '  -3.40%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '29.59'
$ws.Cells.Item(26, 4).NumberFormat = "General"
$ws.Cells.Item(26, 5).Value = '  +0.73%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -0.07%  '

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '10.30'
$ws.Cells.Item(28, 4).NumberFormat = "General"
$ws.Cells.Item(28, 5).Value = '  -3.43%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '37.62'
$ws.Cells.Item(29, 4).NumberFormat = "General"
$ws.Cells.Item(29, 5).Value = '  -3.73%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'Toncoin'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '2.15'
$ws.Cells.Item(30, 4).NumberFormat = "General"
$ws.Cells.Item(30, 5).Value = '  -4.73%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -4.49%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.61'
$ws.Cells.Item(32, 4).NumberFormat = "General"
$ws.Cells.Item(32, 5).Value = '  +0.42%  '

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '2.26'
$ws.Cells.Item(33, 4).NumberFormat = "General"
$ws.Cells.Item(33, 5).Value = '  +2.68%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'WEMIXToken'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '2.80'
$ws.Cells.Item(34, 4).NumberFormat = "General"
$ws.Cells.Item(34, 5).Value = '  -1.51%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'Monero'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '151.96'
$ws.Cells.Item(35, 4).NumberFormat = "General"
$ws.Cells.Item(35, 5).Value = '  +0.45%  '

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.0819'
$ws.Cells.Item(36, 4).NumberFormat = "General"

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.117'
$ws.Cells.Item(37, 4).NumberFormat = "General"
$ws.Cells.Item(37, 5).Value = '  -4.64%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'Stellar'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.121'
$ws.Cells.Item(38, 4).NumberFormat = "General"
$ws.Cells.Item(38, 5).Value = '  -1.37%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'EnergySwap'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '24.43'
$ws.Cells.Item(39, 4).NumberFormat = "General"
$ws.Cells.Item(39, 5).Value = '  +5.35%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '17.01'
$ws.Cells.Item(40, 4).NumberFormat = "General"
$ws.Cells.Item(40, 5).Value = '  +5.04%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '3.56'
$ws.Cells.Item(41, 4).NumberFormat = "General"
$ws.Cells.Item(41, 5).Value = '  -2.53%  '

# Row 42
$ws.Cells.Item(42, 5).Value = '  -4.78%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '3.86'
$ws.Cells.Item(43, 4).NumberFormat = "General"
$ws.Cells.Item(43, 5).Value = '  -5.70%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '2.069.14'
$ws.Cells.Item(44, 5).Value = '  -3.59%  '

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.997'
$ws.Cells.Item(45, 4).NumberFormat = "General"
$ws.Cells.Item(45, 5).Value = '  -0.02%  '

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '88.91'
$ws.Cells.Item(46, 4).NumberFormat = "General"
$ws.Cells.Item(46, 5).Value = '  -4.35%  '

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '9.08'
$ws.Cells.Item(47, 4).NumberFormat = "General"
$ws.Cells.Item(47, 5).Value = '  -4.16%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(48, 4).Value = '2.848.75'
$ws.Cells.Item(48, 5).Value = '  +0.24%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.59'
$ws.Cells.Item(49, 4).NumberFormat = "General"
$ws.Cells.Item(49, 5).Value = '  +3.33%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '106.31'
$ws.Cells.Item(50, 4).NumberFormat = "General"
$ws.Cells.Item(50, 5).Value = '  -2.74%  '

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.191'
$ws.Cells.Item(51, 4).NumberFormat = "General"
$ws.Cells.Item(51, 5).Value = '  -4.95%  '
